$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K2: 50 -> 5 (the Raspberry Pi router is slower than expected, so the
# connectivity-delay score input is lowered)
$ws.Range("K2").Value = 5

# T2: the EXP() latency-decay term no longer divides K2 by 10, making the
# scoring curve easier to satisfy at the new (smaller) delay values
$ws.Range("T2").Formula = "=3*MIN(5,F2+H2+J2+L2+N2+P2+R2)+D2*10+E2*10+10*EXP(-G2/100)+30*(1-EXP(-I2/20))+30*(1-EXP(-K2))+10*M2+15*O2+15*Q2+S2"

# Move the active selection in the frozen bottom-right pane to K3
$ws.Range("K3").Select()
